# Applies the cryptos.xlsx data refresh described by the commit diff.
# Forces column D/B/C string values to remain text (matching the source
# workbook, where these cells are stored as inline strings) even when the
# text looks like a number, by using the classic leading-apostrophe text
# marker and then resetting the cell style so no stray number format is
# left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 "25.227.03"
Set-TextCell 2 5 "  -2.83%  "

# Row 3
Set-TextCell 3 4 "1.548.15"
Set-TextCell 3 5 "  -4.88%  "

# Row 4
Set-TextCell 4 5 "  -0.11%  "

# Row 5
Set-TextCell 5 4 "206.75"
Set-TextCell 5 5 "  -3.51%  "

# Row 6
Set-TextCell 6 5 "  -0.08%  "

# Row 7
Set-TextCell 7 4 "0.477"
Set-TextCell 7 5 "  -5.40%  "

# Row 8
Set-TextCell 8 4 "0.0608"
Set-TextCell 8 5 "  -1.63%  "

# Row 9
Set-TextCell 9 4 "0.242"
Set-TextCell 9 5 "  -3.25%  "

# Row 10
Set-TextCell 10 4 "17.73"
Set-TextCell 10 5 "  -4.08%  "

# Row 11
Set-TextCell 11 5 "  -1.23%  "

# Row 12
Set-TextCell 12 4 "1.764.18"
Set-TextCell 12 5 "  -4.82%  "

# Row 13
Set-TextCell 13 4 "1.558.52"
Set-TextCell 13 5 "  -4.82%  "

# Row 14
Set-TextCell 14 5 "  -4.72%  "

# Row 15
Set-TextCell 15 5 "  -4.62%  "

# Row 16
Set-TextCell 16 4 "25.187.01"
Set-TextCell 16 5 "  -3.00%  "

# Row 17
Set-TextCell 17 5 "  -4.54%  "

# Row 18
Set-TextCell 18 4 "58.49"
Set-TextCell 18 5 "  -4.74%  "

# Row 19
Set-TextCell 19 4 "1.01"
Set-TextCell 19 5 "  -0.07%  "

# Row 20
Set-TextCell 20 4 "185.56"
Set-TextCell 20 5 "  -3.74%  "

# Row 21
Set-TextCell 21 4 "4.09"
Set-TextCell 21 5 "  -3.73%  "

# Row 22
Set-TextCell 22 4 "9.22"
Set-TextCell 22 5 "  -3.79%  "

# Row 23
Set-TextCell 23 4 "5.81"
Set-TextCell 23 5 "  -4.29%  "

# Row 24
Set-TextCell 24 5 "  -0.08%  "

# Row 25
Set-TextCell 25 5 "  -3.96%  "

# Row 26
Set-TextCell 26 4 "139.07"
Set-TextCell 26 5 "  -3.54%  "

# Row 27
Set-TextCell 27 5 "  -5.45%  "

# Row 28
Set-TextCell 28 5 "  -2.93%  "

# Row 29
Set-TextCell 29 4 "6.37"
Set-TextCell 29 5 "  -5.42%  "

# Row 30
Set-TextCell 30 5 "  -6.62%  "

# Row 31
Set-TextCell 31 5 "  -4.70%  "

# Row 32
Set-TextCell 32 5 "  -3.86%  "

# Row 33
Set-TextCell 33 4 "2.96"
Set-TextCell 33 5 "  -5.02%  "

# Row 34
Set-TextCell 34 5 "  -3.57%  "

# Row 35
Set-TextCell 35 4 "2.31"
Set-TextCell 35 5 "  -3.92%  "

# Row 36
Set-TextCell 36 4 "1.083.32"
Set-TextCell 36 5 "  -3.76%  "

# Row 37
Set-TextCell 37 5 "  -0.50%  "

# Row 38
Set-TextCell 38 5 "  -2.85%  "

# Row 39
Set-TextCell 39 4 "0.492"
Set-TextCell 39 5 "  -5.70%  "

# Row 40
Set-TextCell 40 2 "TrustWalletToken"
Set-TextCell 40 3 "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell 40 4 "0.804"
Set-TextCell 40 5 "  +4.94%  "

# Row 41
Set-TextCell 41 2 "MXToken"
Set-TextCell 41 3 "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell 41 4 "2.24"
Set-TextCell 41 5 "  -7.79%  "

# Row 42
Set-TextCell 42 2 "ARBITRUM"
Set-TextCell 42 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell 42 4 "0.759"
Set-TextCell 42 5 "  -10.91%  "

# Row 43
Set-TextCell 43 5 "  -5.68%  "

# Row 44
Set-TextCell 44 4 "5.03"
Set-TextCell 44 5 "  -2.58%  "

# Row 45
Set-TextCell 45 4 "1.679.81"
Set-TextCell 45 5 "  -4.78%  "

# Row 46
Set-TextCell 46 4 "0.0₆0111"
Set-TextCell 46 5 "  -3.20%  "

# Row 47
Set-TextCell 47 5 "  -1.83%  "

# Row 48
Set-TextCell 48 4 "52.22"
Set-TextCell 48 5 "  -4.16%  "

# Row 49
Set-TextCell 49 4 "0.0502"
Set-TextCell 49 5 "  -5.83%  "

# Row 50
Set-TextCell 50 5 "  -0.27%  "

# Row 51
Set-TextCell 51 4 "0.403"
Set-TextCell 51 5 "  -2.23%  "
